$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 98.912777
$ws.Range("H2").Value = 296.738331
$ws.Range("I2").Value = 0.8120825131376513
$ws.Range("J2").Value = 0.8120825131376513
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 30.84053233333333
$ws.Range("N2").Value = 92.521597
$ws.Range("O2").Value = 0.2792762306509579
$ws.Range("P2").Value = 0.2792762306509579
$ws.Range("Q2").Value = 3050.52269724829
$ws.Range("R2").Value = 27454.70427523461
$ws.Range("S2").Value = 0.2267953432466403
$ws.Range("T2").Value = 0.2267953432466403

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 98.912777
$ws.Range("H3").Value = 296.738331
$ws.Range("I3").Value = 0.8120825131376513
$ws.Range("J3").Value = 0.8120825131376513
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 69.52746833333333
$ws.Range("N3").Value = 208.582405
$ws.Range("O3").Value = 0.6296055163046042
$ws.Range("P3").Value = 0.6296055163046042
$ws.Range("Q3").Value = 6877.154970629562
$ws.Range("R3").Value = 61894.39473566606
$ws.Range("S3").Value = 0.5112916299659714
$ws.Range("T3").Value = 0.5112916299659714

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 98.912777
$ws.Range("H4").Value = 296.738331
$ws.Range("I4").Value = 0.8120825131376513
$ws.Range("J4").Value = 0.8120825131376513
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 10.06220766666667
$ws.Range("N4").Value = 30.186623
$ws.Range("O4").Value = 0.09111825304443795
$ws.Range("P4").Value = 0.09111825304443795
$ws.Range("Q4").Value = 995.2809030606903
$ws.Range("R4").Value = 8957.528127546213
$ws.Range("S4").Value = 0.07399553992503963
$ws.Range("T4").Value = 0.07399553992503963

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 17.04862266666667
$ws.Range("H5").Value = 51.14586800000001
$ws.Range("I5").Value = 0.1399706767982279
$ws.Range("J5").Value = 0.1399706767982279
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 30.84053233333333
$ws.Range("N5").Value = 92.521597
$ws.Range("O5").Value = 0.2792762306509579
$ws.Range("P5").Value = 0.2792762306509579
$ws.Range("Q5").Value = 525.788598590133
$ws.Range("R5").Value = 4732.097387311197
$ws.Range("S5").Value = 0.03909048301787259
$ws.Range("T5").Value = 0.03909048301787259

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 17.04862266666667
$ws.Range("H6").Value = 51.14586800000001
$ws.Range("I6").Value = 0.1399706767982279
$ws.Range("J6").Value = 0.1399706767982279
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 69.52746833333333
$ws.Range("N6").Value = 208.582405
$ws.Range("O6").Value = 0.6296055163046042
$ws.Range("P6").Value = 0.6296055163046042
$ws.Range("Q6").Value = 1185.347572583616
$ws.Range("R6").Value = 10668.12815325254
$ws.Range("S6").Value = 0.08812631023305319
$ws.Range("T6").Value = 0.08812631023305319

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 17.04862266666667
$ws.Range("H7").Value = 51.14586800000001
$ws.Range("I7").Value = 0.1399706767982279
$ws.Range("J7").Value = 0.1399706767982279
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 10.06220766666667
$ws.Range("N7").Value = 30.186623
$ws.Range("O7").Value = 0.09111825304443795
$ws.Range("P7").Value = 0.09111825304443795
$ws.Range("Q7").Value = 171.5467817026405
$ws.Range("R7").Value = 1543.921035323764
$ws.Range("S7").Value = 0.01275388354730217
$ws.Range("T7").Value = 0.01275388354730217

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5.839988000000001
$ws.Range("H8").Value = 17.519964
$ws.Range("I8").Value = 0.0479468100641207
$ws.Range("J8").Value = 0.04794681006412069
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 30.84053233333333
$ws.Range("N8").Value = 92.521597
$ws.Range("O8").Value = 0.2792762306509579
$ws.Range("P8").Value = 0.2792762306509579
$ws.Range("Q8").Value = 180.1083387402787
$ws.Range("R8").Value = 1620.975048662508
$ws.Range("S8").Value = 0.01339040438644504
$ws.Range("T8").Value = 0.01339040438644504

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5.839988000000001
$ws.Range("H9").Value = 17.519964
$ws.Range("I9").Value = 0.0479468100641207
$ws.Range("J9").Value = 0.04794681006412069
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 69.52746833333333
$ws.Range("N9").Value = 208.582405
$ws.Range("O9").Value = 0.6296055163046042
$ws.Range("P9").Value = 0.6296055163046042
$ws.Range("Q9").Value = 406.0395807370467
$ws.Range("R9").Value = 3654.35622663342
$ws.Range("S9").Value = 0.0301875761055795
$ws.Range("T9").Value = 0.0301875761055795

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 5.839988000000001
$ws.Range("H10").Value = 17.519964
$ws.Range("I10").Value = 0.0479468100641207
$ws.Range("J10").Value = 0.0479468100641207
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 10.06220766666667
$ws.Range("N10").Value = 30.186623
$ws.Range("O10").Value = 0.09111825304443795
$ws.Range("P10").Value = 0.09111825304443795
$ws.Range("Q10").Value = 58.76317202684134
$ws.Range("R10").Value = 528.8685482415721
$ws.Range("S10").Value = 0.004368829572096155
$ws.Range("T10").Value = 0.004368829572096154
